$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for this product/market combo.
# It becomes the new row 11 (the block is ordered by date), pushing the
# previously-existing rows 11-14 down to rows 12-15.
$ws.Rows.Item(11).Insert()

# Fill in the data for the newly inserted row 11.
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 45044
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = "Otros"
$ws.Cells.Item(11, 9).Value = 100107011
$ws.Cells.Item(11, 10).Value = "Tuna"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 60
$ws.Cells.Item(11, 14).Value = 15000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 15000
$ws.Cells.Item(11, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(11, 18).Value = "Región Metropolitana"
$ws.Cells.Item(11, 19).Value = 833
$ws.Cells.Item(11, 20).Value = 18
